$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSheet")

for ($r = 12; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "Abhi_0"
    $ws.Cells.Item($r, 2).Value = "Abhi_1"
    $ws.Cells.Item($r, 3).Value = "Abhi_2"
}
